$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 01:52"

# Country list reshuffled upstream: Jamaica now sorts ahead of Monaco, Congo ahead of
# Barbados, and Eritrea ahead of Libia. Rewrite the country-name column for the affected rows.
$ws.Range("A134").Value = "Jamaica"
$ws.Range("A135").Value = "Monaco"
$ws.Range("A136").Value = "Aruba"
$ws.Range("A137").Value = "Guayana Francesa"
$ws.Range("A138").Value = "Etiopia"
$ws.Range("A139").Value = "Liechtenstein"
$ws.Range("A140").Value = "Togo"
$ws.Range("A141").Value = "Congo"
$ws.Range("A142").Value = "Barbados"
$ws.Range("A143").Value = "Birmania"
$ws.Range("A144").Value = "Somalia"
$ws.Range("A161").Value = "Eritrea"
$ws.Range("A162").Value = "Libia"
$ws.Range("A163").Value = "Benin"

# Refresh the numeric columns (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes) with the latest figures.
$ws.Range("B4").Value = 611745
$ws.Range("C4").Value = 24804
$ws.Range("E4").Value = 547095
$ws.Range("F4").Value = 13473
$ws.Range("G4").Value = 2335
$ws.Range("H4").Value = 25975

$ws.Range("B5").Value = 174060
$ws.Range("C5").Value = 3961
$ws.Range("E5").Value = 88301
$ws.Range("G5").Value = 499
$ws.Range("H5").Value = 18255

$ws.Range("B8").Value = 132210
$ws.Range("C8").Value = 2138
$ws.Range("E8").Value = 60515
$ws.Range("G8").Value = 301
$ws.Range("H8").Value = 3495

$ws.Range("B96").Value = 492
$ws.Range("C96").Value = 9
$ws.Range("D96").Value = 260
$ws.Range("E96").Value = 224
$ws.Range("F96").Value = 14

$ws.Range("B134").Value = 105
$ws.Range("C134").Value = 32
$ws.Range("D134").Value = 19
$ws.Range("E134").Value = 82
$ws.Range("F134").Value = 0
$ws.Range("H134").Value = 4

$ws.Range("B135").Value = 93
$ws.Range("D135").Value = 6
$ws.Range("E135").Value = 86
$ws.Range("F135").Value = 5
$ws.Range("H135").Value = 1

$ws.Range("B136").Value = 92
$ws.Range("D136").Value = 32
$ws.Range("E136").Value = 60

$ws.Range("B137").Value = 86
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 51
$ws.Range("E137").Value = 35
$ws.Range("F137").Value = 1
$ws.Range("H137").Value = 0

$ws.Range("B138").Value = 82
$ws.Range("C138").Value = 8
$ws.Range("D138").Value = 14
$ws.Range("E138").Value = 65
$ws.Range("H138").Value = 3

$ws.Range("B139").Value = 79
$ws.Range("D139").Value = 55
$ws.Range("E139").Value = 23
$ws.Range("H139").Value = 1

$ws.Range("B140").Value = 77
$ws.Range("D140").Value = 32
$ws.Range("E140").Value = 42
$ws.Range("H140").Value = 3

$ws.Range("B141").Value = 74
$ws.Range("C141").Value = 14
$ws.Range("D141").Value = 10
$ws.Range("E141").Value = 59
$ws.Range("F141").Value = 0
$ws.Range("H141").Value = 5

$ws.Range("B142").Value = 73
$ws.Range("D142").Value = 15
$ws.Range("E142").Value = 53
$ws.Range("F142").Value = 4
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 5

$ws.Range("B143").Value = 63
$ws.Range("C143").Value = 1
$ws.Range("E143").Value = 57
$ws.Range("F143").Value = 0
$ws.Range("H143").Value = 4

$ws.Range("D144").Value = 2
$ws.Range("E144").Value = 56
$ws.Range("F144").Value = 2
$ws.Range("H144").Value = 2

$ws.Range("C161").Value = 1
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 35
$ws.Range("H161").Value = 0

$ws.Range("C162").Value = 9
$ws.Range("D162").Value = 9
$ws.Range("E162").Value = 25

$ws.Range("B163").Value = 35
$ws.Range("D163").Value = 18
$ws.Range("E163").Value = 16
$ws.Range("H163").Value = 1
